# PlanProjekta.xlsx edit — "Add ERD, ConceptualModel, Fix UCD"
#
# Summary of changes applied:
#  - B1 header "Trajanje" -> "Trajanje(dani)"
#  - A3 "Procjena resursa" -> "Procjena resursa, ciljeva"
#  - A4 "Plan intervjurianja" -> "Istraživanje tržišta"
#  - A6/A7 swapped: A6 becomes "Plan projekta", A7 becomes "Analiza izvedivosti"
#  - Column B widened
#  - Active selection moved to A8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates -------------------------------------------------
# Order matters: new shared-string entries are appended in first-use order,
# and the target file order is "Procjena resursa, ciljeva", "Istraživanje
# tržišta", then "Trajanje(dani)" — so write A3/A4 before B1.
$ws.Range("A3").Value = "Procjena resursa, ciljeva"
$ws.Range("A4").Value = "Istraživanje tržišta"
$ws.Range("A6").Value = "Plan projekta"
$ws.Range("A7").Value = "Analiza izvedivosti"
$ws.Range("B1").Value = "Trajanje(dani)"

# --- Column B width ---------------------------------------------------
# Target stored width is 21.5546875 chars; the COM ColumnWidth setter here
# quantizes to 1/6-character steps, so 20.6667 is the closest input that
# lands on the nearest reachable stored width (21.5).
$ws.Columns.Item(2).ColumnWidth = 20.6667

# --- Selection / active cell -------------------------------------------
$ws.Activate()
$ws.Range("A8").Select()
